$wb = $excel.ActiveWorkbook

# --- Netherlands: copy of Portugal layout (21-row full product list) ---
$portugal = $wb.Worksheets.Item("Portugal")
$portugal.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$netherlands = $wb.Worksheets.Item($wb.Worksheets.Count)
$netherlands.Name = "Netherlands"
$netherlands.Range("B4").Value = "NGC-3144/T2199"
$netherlands.Range("B2").Value = "Netherlands Market"

# --- Austria: copy of Slovakia layout (shorter product list) ---
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$austria = $wb.Worksheets.Item($wb.Worksheets.Count)
$austria.Name = "Austria"
$austria.Range("B4").Value = "NGC-3817/T2306"
$austria.Range("B2").Value = "Austria Market"

# --- Denmark: copy of Slovakia layout plus an extra product row ---
$slovakia2 = $wb.Worksheets.Item("Slovakia")
$slovakia2.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"

$denmark.Rows.Item(18).Insert()
$denmark.Range("A17").Copy()
$denmark.Range("A18").PasteSpecial(-4122)
$denmark.Range("A18").Value = "MZXSDR240"

$denmark.Range("B4").Value = "NGC-2913/T2798"
$denmark.Range("B2").Value = "Denmark Market"

# Fix up selections so each new sheet's cursor sits on B4 (matches source sheets)
$netherlands.Range("B4").Select()
$austria.Range("B4").Select()
$denmark.Range("B4").Select()

# Netherlands should end up the active tab (mirrors upstream workbook state)
$netherlands.Activate()
$netherlands.Range("B4").Select()

Write-Output "done"
